# Update res_bus vm_pu values for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order matching the data table below: B,C,D,E,F,I,J,K,L,M,N
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# Each row below corresponds to worksheet rows 2..25 (A2..A25 = 0..23)
$data = @(
    @(1.02,1.025203800266516,1.053124818602783,1.025592506232896,1.057016744503547,1.041008870720184,1.03037413121906,1.055871652524516,1.028417625384855,1.059752868684721,1.014162600421638),
    @(1.02,1.025978396270902,1.053743977164223,1.026245666950282,1.057763787259586,1.041181266677427,1.03078886405342,1.056304290386878,1.028878540275439,1.060313834279636,1.014301233051976),
    @(1.02,1.026480401254202,1.054145152940768,1.026669376181177,1.058248113473574,1.041292008847004,1.031057338907496,1.05658408016659,1.029177168360221,1.06067708706102,1.014390940003419),
    @(1.02,1.026691631299372,1.054313934122284,1.026847757875631,1.058451946893091,1.041338370222849,1.031170231956814,1.056701664886034,1.029302802409559,1.060829861442,1.014428652896054),
    @(1.02,1.026727108673518,1.054342280604775,1.026877723804166,1.058486184413326,1.041346143064532,1.031189188696183,1.056721405550114,1.029323902201883,1.060855516568046,1.014434985048086),
    @(1.02,1.026483222987378,1.05414740770725,1.026671758729224,1.058250836233571,1.04129262909578,1.031058847288555,1.056585651494234,1.029178846733564,1.060679128194906,1.014391443925181),
    @(1.02,1.025465413998526,1.053333953315959,1.025813022107092,1.05726901519106,1.041067299865152,1.030514267530606,1.056017895826827,1.028573313075778,1.059942392504485,1.014209451190849),
    @(1.02,1.023678034762174,1.051904762093037,1.024308104550669,1.055546217330628,1.04066407744947,1.029555586759834,1.055016314703544,1.027509300535004,1.058646324036989,1.013888795978795),
    @(1.02,1.022490682505065,1.050954930557432,1.02331051477969,1.054402726477172,1.040391168057059,1.028917175248709,1.054347934279023,1.026802072597778,1.057783837583779,1.013675080745094),
    @(1.02,1.021977573866034,1.050544372411314,1.022879921610016,1.053908806743609,1.040272034630785,1.028640920740197,1.054058381956929,1.026496354810055,1.057410764801405,1.01358255875893),
    @(1.02,1.02178713805256,1.050391983963421,1.02272018797141,1.053725528504605,1.040227639504115,1.028538336185393,1.053950810063023,1.026382876716472,1.057272249399362,1.013548195186474),
    @(1.02,1.021827980109636,1.050424666739098,1.022754441941318,1.053764833875198,1.040237168897248,1.0285603396017,1.053973885445047,1.026407214549108,1.057301958657505,1.013555566136166),
    @(1.02,1.021961829197315,1.050531773655121,1.022866713726017,1.053893653108634,1.040268367841677,1.028632440476274,1.054049490414216,1.026486973053085,1.057399313831024,1.013579718186939),
    @(1.02,1.022044318711556,1.050597780551633,1.022935915679413,1.053973047546197,1.040287571528832,1.02867686804393,1.054096070626741,1.0265361254266,1.057459305632766,1.013594599502915),
    @(1.02,1.022524757547738,1.050982193414098,1.023339120864422,1.054435532194035,1.040399054356678,1.02893551329747,1.054367148120711,1.026822373091823,1.057808605566475,1.013681221548249),
    @(1.02,1.022826399309535,1.051223521135822,1.023592408997346,1.054725964708444,1.040468727665444,1.029097804276847,1.054537151766154,1.027002068047753,1.058027817831401,1.013735562466451),
    @(1.02,1.023002440394364,1.051364353336133,1.023740279707542,1.054895486572827,1.040509274065831,1.029192483318949,1.05463629847666,1.027106930817275,1.058155718093651,1.01376726032016),
    @(1.02,1.023062482511883,1.051412385253913,1.02379072213038,1.054953309005031,1.040523083564016,1.029224769346155,1.054670102596223,1.027142694703367,1.058199335072867,1.013778068745179),
    @(1.02,1.022794025795965,1.051197621717342,1.023565219906824,1.054694791878551,1.040461261977875,1.029080390175058,1.054518913372527,1.026982783323777,1.058004294541824,1.013729732019227),
    @(1.02,1.021922409666792,1.050500230252781,1.022833646730029,1.053855713914417,1.040259184489331,1.028611207763566,1.054027227153941,1.026463483965912,1.057370643475108,1.013572605917338),
    @(1.02,1.021375291113245,1.050062396585409,1.022374881125353,1.053329226943349,1.040131299551647,1.028316380699173,1.053717973839698,1.026137438263808,1.056972592536416,1.013473833330998),
    @(1.02,1.021665242951769,1.050294438724092,1.022617966694865,1.053608225000698,1.04019917224601,1.028472657869971,1.053881924815465,1.026310237308245,1.057183572994846,1.01352619259939),
    @(1.02,1.022808653689003,1.051209324334281,1.023577505078828,1.054708877175849,1.04046463568777,1.029088258804308,1.0545271545596,1.026991497106949,1.058014923588359,1.013732366542124),
    @(1.02,1.024139377112836,1.052273729296523,1.024696167748393,1.055546217330628,1.040769045417153,1.029803309469629,1.055016314703544,1.027509300535004,1.058646324036989,1.013888795978795)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $row).Value = $values[$j]
    }
}
